$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generator (DistMult) results
$ws.Range("F3").Value = 0.018585
$ws.Range("G3").Value = 123.259453
$ws.Range("H3").Value = 0.014372

# Discriminator (TransE) results
$ws.Range("K3").Value = 0.020653
$ws.Range("L3").Value = 123.334343
$ws.Range("M3").Value = 0.017398

# Adversarial Training results
$ws.Range("O3").Value = 0.020921
$ws.Range("P3").Value = 123.322998
$ws.Range("Q3").Value = 0.015885

# Update the view selection to match the saved window state
$ws.Range("R3").Select()
